$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Use Cases" overview table (3 columns: Use Case / Aktor / Kapitel)
#    Fill in the previously-empty second row and append nine more rows that
#    list each use case together with its actor(s).
# ---------------------------------------------------------------------------
$useCaseTable = $d.Tables.Item(5)

# the row that already exists (currently three empty cells)
$useCaseTable.Cell(2, 1).Range.Text = "Choose Sport"
$useCaseTable.Cell(2, 2).Range.Text = "Trainer"

$useCaseRows = @(
    @("Create Board", "Trainer"),
    @("Save Board", "Trainer"),
    @("Load Board", "Trainer"),
    @("Delete Board", "Trainer"),
    @("Change Field", "Trainer"),
    @("Draw Item", "Trainer, Jugend und Sport"),
    @("Change Item", "Trainer"),
    @("Delete item", "Trainer"),
    @("Add Attributes", "Trainer")
)

foreach ($pair in $useCaseRows) {
    $newRow = $useCaseTable.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $pair[0]
    $newRow.Cells.Item(2).Range.Text = $pair[1]
}

# ---------------------------------------------------------------------------
# 2) "Aktoren" table (2 columns: Aktor / Beschreibung)
#    Fill in the previously-empty second row and append two more rows.
# ---------------------------------------------------------------------------
$aktorenTable = $d.Tables.Item(6)

$aktorenTable.Cell(2, 1).Range.Text = "Trainer"
$aktorenTable.Cell(2, 2).Range.Text = "Primäraktor: Er ist der Endnutzer vom TBE, zeichnet damit Übungen und gestaltet Trainingseinheiten."

$aktorenRows = @(
    @("Jugend und Sport", "Unterstützender Aktor: Liefert die Standart-Notation der einzelnen Sportarten"),
    @("Sportverbände", "Offstage-Aktoren: Sie haben ein Interesse daran, dass die allgemeine Trainingsqualität gewährleistet ist.")
)

foreach ($pair in $aktorenRows) {
    $newRow = $aktorenTable.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $pair[0]
    $newRow.Cells.Item(2).Range.Text = $pair[1]
}

# ---------------------------------------------------------------------------
# 3) Rename the "Use Case 1" heading to "Choose Sport" and drop its TOC
#    bookmark (the heading no longer needs to be a distinct TOC target).
#    Use the bookmark's own range to retarget only the heading occurrence
#    (the document also contains an unrelated "Use Case 1" entry inside the
#    table of contents that must stay untouched).
# ---------------------------------------------------------------------------
$ucBookmark = $d.Bookmarks.Item("_Toc163986095")
$ucBookmark.Range.Text = "Choose Sport"
$d.Bookmarks.Item("_Toc163986095").Delete()
